# Update "想去人数" (F column) values for events that changed
# between the two crawls, on both the "展览" sheet and the
# "全部类型" sheet (which mirrors 展览 plus an extra row).

$wb = $excel.ActiveWorkbook

# Sheet "展览" (first sheet) - row => new F value
$exhibitionUpdates = @{
    2  = 6868
    3  = 16
    4  = 443
    9  = 107
    12 = 41
    13 = 185
    14 = 433
    15 = 17
    16 = 1791
    17 = 32
    18 = 3482
    20 = 238
    21 = 20
    22 = 2116
    23 = 196
    24 = 7
    25 = 30
    26 = 5
}

$wsExhibition = $wb.Worksheets.Item("展览")
foreach ($row in $exhibitionUpdates.Keys) {
    $wsExhibition.Cells.Item($row, 6).Value = $exhibitionUpdates[$row]
}

# Sheet "全部类型" (fourth sheet) - same events, shifted by +1 row
# starting at row 10 because of an extra entry at row 7.
$allTypesUpdates = @{
    2  = 6868
    3  = 16
    4  = 443
    10 = 107
    13 = 41
    14 = 185
    15 = 433
    16 = 17
    17 = 1791
    18 = 32
    19 = 3482
    21 = 238
    22 = 20
    23 = 2116
    24 = 196
    25 = 7
    26 = 30
    27 = 5
}

$wsAllTypes = $wb.Worksheets.Item("全部类型")
foreach ($row in $allTypesUpdates.Keys) {
    $wsAllTypes.Cells.Item($row, 6).Value = $allTypesUpdates[$row]
}
